$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (avoids Excel
# silently re-typing numeric-looking strings like "246.85" as numbers),
# then restores the cell's original style so no stray style/numFmt
# is left behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}

Set-TextValue 2 4 '42.217.23'
$ws.Cells.Item(2, 5).Value = '  -1.57%  '
Set-TextValue 3 4 '2.245.75'
$ws.Cells.Item(3, 5).Value = '  -2.01%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
Set-TextValue 5 4 '246.85'
$ws.Cells.Item(5, 5).Value = '  -2.06%  '
Set-TextValue 6 4 '0.628'
$ws.Cells.Item(6, 5).Value = '  -2.74%  '
Set-TextValue 7 4 '74.36'
$ws.Cells.Item(7, 5).Value = '  -1.26%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
Set-TextValue 9 4 '0.615'
$ws.Cells.Item(9, 5).Value = '  -5.21%  '
Set-TextValue 10 4 '41.04'
$ws.Cells.Item(10, 5).Value = '  +5.20%  '
Set-TextValue 11 4 '0.0941'
$ws.Cells.Item(11, 5).Value = '  -4.40%  '
Set-TextValue 12 4 '7.09'
$ws.Cells.Item(12, 5).Value = '  -5.74%  '
Set-TextValue 13 4 '0.103'
$ws.Cells.Item(13, 5).Value = '  -3.89%  '
Set-TextValue 14 4 '2.576.29'
$ws.Cells.Item(14, 5).Value = '  -2.40%  '
Set-TextValue 15 4 '14.49'
$ws.Cells.Item(15, 5).Value = '  -4.28%  '
Set-TextValue 16 4 '0.852'
$ws.Cells.Item(16, 5).Value = '  -2.37%  '
Set-TextValue 17 4 '2.229.04'
Set-TextValue 18 4 '42.066.62'
$ws.Cells.Item(18, 5).Value = '  -1.74%  '
Set-TextValue 19 4 '0.0₃0973'
$ws.Cells.Item(19, 5).Value = '  -2.81%  '
Set-TextValue 20 4 '71.68'
$ws.Cells.Item(20, 5).Value = '  -0.84%  '
Set-TextValue 21 4 '6.08'
$ws.Cells.Item(21, 5).Value = '  -2.39%  '
Set-TextValue 22 4 '2.33'
$ws.Cells.Item(22, 5).Value = '  +7.71%  '
Set-TextValue 23 4 '230.85'
$ws.Cells.Item(23, 5).Value = '  -2.06%  '
$ws.Cells.Item(24, 5).Value = '  +0.06%  '
Set-TextValue 25 4 '11.08'
$ws.Cells.Item(25, 5).Value = '  -2.49%  '
Set-TextValue 26 4 '3.56'
$ws.Cells.Item(26, 5).Value = '  -8.04%  '
$ws.Cells.Item(27, 5).Value = '  -4.41%  '
Set-TextValue 28 4 '7.44'
$ws.Cells.Item(28, 5).Value = '  +19.88%  '
Set-TextValue 29 4 '2.15'
$ws.Cells.Item(29, 5).Value = '  -1.06%  '
Set-TextValue 30 4 '169.49'
$ws.Cells.Item(30, 5).Value = '  +1.23%  '
Set-TextValue 31 4 '20.67'
$ws.Cells.Item(31, 5).Value = '  -1.67%  '
Set-TextValue 32 4 '0.0827'
$ws.Cells.Item(32, 5).Value = '  -4.13%  '
$ws.Cells.Item(33, 5).Value = '  -5.88%  '
Set-TextValue 34 4 '30.13'
$ws.Cells.Item(34, 5).Value = '  -4.42%  '
$ws.Cells.Item(35, 5).Value = '  -2.69%  '
Set-TextValue 36 4 '4.54'
$ws.Cells.Item(36, 5).Value = '  -2.52%  '
$ws.Cells.Item(37, 5).Value = '  +0.68%  '
$ws.Cells.Item(38, 5).Value = '  -2.14%  '
Set-TextValue 39 4 '13.35'
$ws.Cells.Item(39, 5).Value = '  -1.83%  '
$ws.Cells.Item(40, 5).Value = '  -5.61%  '
$ws.Cells.Item(41, 5).Value = '  -3.46%  '
Set-TextValue 42 4 '109.66'
$ws.Cells.Item(42, 5).Value = '  +3.75%  '
Set-TextValue 43 4 '0.203'
$ws.Cells.Item(43, 5).Value = '  -3.67%  '
Set-TextValue 44 4 '60.94'
$ws.Cells.Item(44, 5).Value = '  -0.14%  '
Set-TextValue 45 4 '8.64'
$ws.Cells.Item(45, 5).Value = '  -5.18%  '
Set-TextValue 46 4 '0.100'
$ws.Cells.Item(46, 5).Value = '  -1.41%  '
Set-TextValue 47 4 '0.998'
$ws.Cells.Item(47, 5).Value = '  -0.21%  '
Set-TextValue 48 4 '1.12'
$ws.Cells.Item(48, 5).Value = '  -4.46%  '
$ws.Cells.Item(49, 5).Value = '  -1.81%  '
Set-TextValue 50 4 '2.24'
$ws.Cells.Item(50, 5).Value = '  -1.41%  '
$ws.Cells.Item(51, 2).Value = 'WOONetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue 51 4 '0.426'
$ws.Cells.Item(51, 5).Value = '  +14.60%  '
